# "added new scripts for new offers" - refresh the DigiZag 30-day report with the
# latest export: 3 brand-new order rows appear at the top (pushing every other
# row down by one date-rank), the running Grand Total/Order Count updates, and
# the 3 rows that used to be empty at the bottom of the table (47-49) are now
# filled in with what had been the tail of the previous export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-49: Row, Discount Code (A), Order Name (B),
# Date serial (C), Order Count (D), Net Sales (E)
$rows = @(
  @(2, "Grand Total", $null, $null, 47.0, 16796.48),
  @(3, "DZZ2", "MB494715", 45935.0, 1.0, 267.3),
  @(4, "DZZ2", "MB494491", 45935.0, 1.0, 351.9),
  @(5, "DZZ2", "MB494422", 45934.0, 1.0, 267.3),
  @(6, "DZCT", "MB494311", 45934.0, 1.0, 534.6),
  @(7, "DZCT", "MB493162", 45933.0, 1.0, 424.8),
  @(8, "DZ65", "MB493466", 45933.0, 1.0, 246.6),
  @(9, "DZ65", "MB492884", 45932.0, 1.0, 424.8),
  @(10, "DZZ2", "MB491923", 45932.0, 1.0, 267.3),
  @(11, "DZCT", "MB492968", 45932.0, 1.0, 437.4),
  @(12, "DZCT", "MB491721", 45931.0, 1.0, 267.3),
  @(13, "DZCT", "MB491606", 45931.0, 1.0, 441.0),
  @(14, "DZZ2", "MB490643", 45930.0, 1.0, 170.1),
  @(15, "DZZ2", "MB489618", 45929.0, 1.0, 267.3),
  @(16, "DZZ2", "MB489648", 45929.0, 1.0, 437.4),
  @(17, "DZCT", "MB489754", 45929.0, 1.0, 522.0),
  @(18, "DZ65", "MB490199", 45929.0, 1.0, 343.8),
  @(19, "DZCT", "MB489048", 45928.0, 1.0, 267.3),
  @(20, "DZZ2", "MB488679", 45928.0, 1.0, 343.8),
  @(21, "DZZ2", "MB489154", 45928.0, 1.0, 267.3),
  @(22, "DZZ2", "MB487512", 45926.0, 1.0, 172.8),
  @(23, "DZZ2", "MB487787", 45926.0, 1.0, 343.8),
  @(24, "DZCT", "MB487316", 45926.0, 1.0, 611.1),
  @(25, "DZCT", "MB487207", 45926.0, 1.0, 343.8),
  @(26, "DZZ2", "MB486945", 45925.0, 1.0, 267.3),
  @(27, "DZCT", "MB486701", 45925.0, 1.0, 351.9),
  @(28, "DZCT", "MB486670", 45925.0, 1.0, 494.1),
  @(29, "DZCT", "MB486223", 45924.0, 1.0, 170.1),
  @(30, "DZCT", "MB485325", 45923.0, 1.0, 343.8),
  @(31, "DZCT", "MB484984", 45922.0, 1.0, 267.3),
  @(32, "DZZ2", "MB484513", 45921.0, 1.0, 178.2),
  @(33, "DZZ2", "MB484651", 45921.0, 1.0, 267.3),
  @(34, "DZCT", "MB484178", 45920.0, 1.0, 437.4),
  @(35, "DZZ2", "MB482319", 45915.0, 1.0, 267.3),
  @(36, "DZZ2", "MB482292", 45915.0, 1.0, 687.6),
  @(37, "DZ77", "MB482234", 45915.0, 1.0, 437.4),
  @(38, "DZZ2", "MB481818", 45914.0, 1.0, 267.3),
  @(39, "DZZ2", "MB481156", 45912.0, 1.0, 343.8),
  @(40, "DZZ2", "MB481115", 45912.0, 1.0, 343.8),
  @(41, "DZCT", "MB480856", 45911.0, 1.0, 1081.8),
  @(42, "DZ65", "MB480838", 45911.0, 1.0, 398.7),
  @(43, "DZZ2", "MB480981", 45911.0, 1.0, 267.3),
  @(44, "DZCT", "MB480389", 45909.0, 1.0, 267.3),
  @(45, "DZ77", "MB479700", 45907.0, 1.0, 267.3),
  @(46, "DZZ2", "MB479337", 45906.0, 1.0, 267.3),
  @(47, "DZ65", "MB479458", 45906.0, 1.0, 563.18),
  @(48, "DZZ2", "MB479273", 45906.0, 1.0, 267.3),
  @(49, "DZ65", "MB479499", 45906.0, 1.0, 270.9)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
